# Updates files with strain names:
#  - Column B (harvester) changes from "S.GISH" to "H.Brown" for all data rows (2-25)
#  - Column F (strain) gets new/updated values for several rows
#  - Selection moves to F21:F22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: harvester -> "H.Brown" for every data row 2..25
$ws.Range("B2:B25").Value = "H.Brown"

# Column F: strain values for specific rows
$ws.Range("F3").Value = "KN99alpha"
$ws.Range("F4").Value = "KN99alpha"

$ws.Range("F6").Value = "TDY2258"
$ws.Range("F7").Value = "TDY2258"

$ws.Range("F9").Value = "TYS2271"
$ws.Range("F10").Value = "TYS2271"

$ws.Range("F20").Value = "TDY1984"
$ws.Range("F21").Value = "TDY1984"
$ws.Range("F22").Value = "TDY1984"

# Update the active selection to match the saved view state
$ws.Range("F21:F22").Select()
